$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new row at position 14 (shifts old rows 14-18 down to 15-19,
# carrying their values / formats / merges with them).
# ------------------------------------------------------------------
$ws.Rows(14).Insert()

# ------------------------------------------------------------------
# The freshly inserted row 14 comes back blank / unformatted, so
# restore its formatting by copying it from the row above (row 13),
# which already carries the exact same layout used by every other
# item row in the table.
# ------------------------------------------------------------------
$ws.Range("A13:Q13").Copy()
$ws.Range("A14:Q14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-create the merges for the new row 14 (same pattern as all other
# item rows: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A14:B14").Merge()
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("L14:M14").Merge()
$ws.Range("N14:O14").Merge()

# Row heights: row 14 takes the "odd" 25.5 height (same as row 16/17),
# matching the final layout.
$ws.Rows(14).RowHeight = 25.5
$ws.Rows(15).RowHeight = 24.75
$ws.Rows(16).RowHeight = 25.5
$ws.Rows(17).RowHeight = 25.5
$ws.Rows(18).RowHeight = 24.75
$ws.Rows(19).RowHeight = 16.5

# ------------------------------------------------------------------
# Populate the new row 14 with the new item: WATER FOR INJECTION.
# The source data stores every figure (balance, order-limit, price,
# selling price, transaction count) as TEXT even though a couple of
# those columns (L, P) carry a numeric-looking display format - so
# those two need a temporary "@" (text) number format to stop Excel
# from re-interpreting the literal as a number on entry.
# ------------------------------------------------------------------
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "WATER FOR INJECTION AMP. 5 ML"
$ws.Range("H14").Value = "8274:0"

$fmtL14 = $ws.Range("L14").NumberFormat
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "1"
$ws.Range("L14").NumberFormat = $fmtL14

$ws.Range("N14").Value = "2.00"

$fmtP14 = $ws.Range("P14").NumberFormat
$ws.Range("P14").NumberFormat = "@"
$ws.Range("P14").Value = "18.0000"
$ws.Range("P14").NumberFormat = $fmtP14

$ws.Range("Q14").Value = "9:0"

# ------------------------------------------------------------------
# Item numbers for the rows that followed shift up by one extra slot:
# old row 14 (ZYRTEC, #8) is now row 15 but keeps number 8; old row 15
# (#9) is now row 16 keeping 9; old row 16 (#10) is now row 17 and
# becomes item #11 (an extra item was inserted ahead of it).
# ------------------------------------------------------------------
$ws.Range("A17").Value = 11

# ------------------------------------------------------------------
# Update the totals row (old row 17 -> now row 18): the selling-price
# sum grows by the new item's selling price (18.0000).
# ------------------------------------------------------------------
$ws.Range("P18").Value = 346.89999999999998

# ------------------------------------------------------------------
# Update the generation timestamp in the footer row (now row 19).
# ------------------------------------------------------------------
$ws.Range("A19").Value = "Saturday, 13 September, 2025 11:08 AM"
